# ===========================================================================
# Report.xlsx: split off a new "Остальные" sheet that keeps the old 250000-
# sample rows, add a PCA+SVM row to the main (5000-sample) sheet, and retitle
# both sheets/tweak their layout.
# ===========================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Create the second sheet right after the first one. Populate it (header +
#    the two rows that used to live at A5 and A7 on sheet1) *before* touching
#    sheet1's own data, via Copy / PasteSpecial so values, shared strings and
#    cell styles all transfer verbatim.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Остальные"

$ws1.Range("A1:G1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4122)
$ws1.Range("A1:G1").Copy()
$ws2.Range("A1:G1").PasteSpecial(-4163)

$ws1.Range("A2:G2").Copy()
$ws2.Range("A2:G3").PasteSpecial(-4122)

$ws1.Range("A5:G5").Copy()
$ws2.Range("A2:G2").PasteSpecial(-4163)

$ws1.Range("A7:G7").Copy()
$ws2.Range("A3:G3").PasteSpecial(-4163)

# Row 4: blank cells that still carry the body formatting.
$ws1.Range("A2:G2").Copy()
$ws2.Range("A4:G4").PasteSpecial(-4122)

$ws2.Range("F1").WrapText = $true

$ws2.Columns.Item(1).ColumnWidth = 26.67
$ws2.Columns.Item(2).ColumnWidth = 26.17
$ws2.Columns.Item(3).ColumnWidth = 15.33
$ws2.Columns.Item(4).ColumnWidth = 17.83
$ws2.Columns.Item(5).ColumnWidth = 23
$ws2.Columns.Item(6).ColumnWidth = 29
$ws2.Columns.Item(7).ColumnWidth = 19.67

$ws2.Rows.Item(1).RowHeight = 75

$ws2.Range("A5").Select()

# ---------------------------------------------------------------------------
# 2) Rewrite sheet1: row 5 takes over the values that used to sit in row 6
#    (kNN leader @ param 12 / 5000 samples), row 6 becomes the brand new
#    RandomizedPCA+SVM method, and the old row 7 (now preserved on sheet 2)
#    is removed.
# ---------------------------------------------------------------------------
$ws1.Range("A6:G6").Copy()
$ws1.Range("A5:G5").PasteSpecial(-4122)
$ws1.Range("A6:G6").Copy()
$ws1.Range("A5:G5").PasteSpecial(-4163)

$ws1.Range("A6").Value = "RandomizedPCA+SVM"
$ws1.Range("B6").Value = "PCA(27, true), SVM(Linear, C=1)"
$ws1.Range("C6").Value = 5000
$ws1.Range("D6").Value = "CV, 5"
$ws1.Range("E6").Value = 0.74
$ws1.Range("F6").Value = 0.01
$ws1.Range("G6").Value = 6.87

$ws1.Rows.Item(7).Delete()

# ---------------------------------------------------------------------------
# 3) Sheet1 cosmetics: rename, widen columns B/F, taller header row, wrap the
#    confidence-interval header, move the selection.
# ---------------------------------------------------------------------------
$ws1.Name = "Тест по выборке из 5000"

$ws1.Columns.Item(2).ColumnWidth = 29.67
$ws1.Columns.Item(6).ColumnWidth = 28.5

$ws1.Rows.Item(1).RowHeight = 28.5

$ws1.Range("F1").WrapText = $true

$ws1.Range("B6").Select()
